$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E7").Value = 65
$ws.Range("D7").Select()
